# Generate Report for Handoff
# This script updates the localization-status report after a new handoff run:
#  - Overview sheet: bump the "Latest HO Xliff Generate Date" for the files
#    that were just re-handed-off
#  - zh-cn / de-de sheets: the four files that were still queued with Priority
#    "low" are now generated with Priority "ht", and their "Latest Handoff
#    Datetime" is refreshed to the new generation timestamp.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: rows 4-7 (32bb56e3, 36eb426b, 3c3fb9b2, 6764dcbf) all shared
# the same "Latest HO Xliff Generate Date" (column G) timestamp, which is bumped
# for the new handoff generation run.
foreach ($row in 4..7) {
    $overview.Range("G$row").Value = "2016-08-16 20:28:35"
}

# zh-cn sheet: rows 4-7 correspond to
#   32bb56e3-c21d-4294-8cc3-96b6c7d5c571.md
#   36eb426b-af30-43d6-ad34-d88496fdf7be.md
#   3c3fb9b2-e5f7-411f-80fe-a8331c525725.md
#   6764dcbf-1d6a-428e-bdcc-5d0490ba27ae.md
# Column E = Priority, Column H = Latest Handoff Datetime
foreach ($row in 4..7) {
    $zhcn.Range("E$row").Value = "ht"
    $zhcn.Range("H$row").Value = "2016-08-16 20:28:30"
}

# de-de sheet: same four rows / files
foreach ($row in 4..7) {
    $dede.Range("E$row").Value = "ht"
    $dede.Range("H$row").Value = "2016-08-16 20:28:35"
}
